# Update significance-annotated p-value strings in column F to remove the
# space between the numeric value and the asterisk(s), across the three
# result sheets (rich, even, invSim).

$wb = $excel.ActiveWorkbook

$wsRich = $wb.Worksheets.Item("rich")
$wsRich.Range("F12").Value = "0.017*"
$wsRich.Range("F13").Value = "0.034*"

$wsEven = $wb.Worksheets.Item("even")
$wsEven.Range("F2").Value = "0.003**"
$wsEven.Range("F11").Value = "0.047*"

$wsInvSim = $wb.Worksheets.Item("invSim")
$wsInvSim.Range("F13").Value = "0.02*"
$wsInvSim.Range("F14").Value = "0.028*"
